$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sponsored`n2 BHK Projects in Thane West | Starts at ₹93 Lacs* by Runwal`nlandsend.runwal.com`nhttp://landsend.runwal.com › projects › thane`nTake Advantage of the Umbrella Offer: 2 BHK Flats Starts at ₹93L* at Lands End by Runwal. Book Your Dream Home at Runwal Lands End And Avail Pay 10% Now & 90% On...`nView Location · View Gallery · Lands End by Runwal · Project Overview · Project Highlights"
$ws.Range("B2").Value = "Realatte Ventures LLP"

$ws.Range("A3").Value = "Sponsored`n1,2,3 BHK by Lodha® in Thane - 1 BHK by Lodha® in Thane`nLodha Group`nhttps://www.lodhagroup.in`nWorld-class 1 BHKs in Thane by India's #1 real estate developer. Building a better life."
$ws.Range("B3").Value = "Macrotech Developers Limited"

$ws.Range("A4").Value = "Sponsored`nGodrej Exquisite, Thane | The Best of Thane`nGodrej Properties`nhttps://www.godrejproperties.com`nNo shared walls between homes | Dual-lobby offers enhanced privacy | Grand Clubhouse. 2&3 BHK Homes starting from ₹1.59 Cr* With Zero Stamp Duty & Registration* Skyscape Gym. Dual-Lobby System. Rooftop Jacuzzi & Spa.`nLocation Details · Property Overview · Podium Amenities · Contact Us · Exclusive Amenities"
$ws.Range("B4").Value = "Madison Communications Private Limited"
